$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2-3: Account holder name and card number
$ws.Range("C2").Value = "Hartmut"

# B3 holds a long digit string that must stay text (not become a number).
# Writing it as a formula that evaluates to a text string, then converting
# that formula to its resulting value in place via PasteSpecial, keeps the
# cell's original style/number format untouched while forcing text storage.
$cardCell = $ws.Range("B3")
$cardCell.Formula = '="2570314725427075"'
$cardCell.Copy()
$cardCell.PasteSpecial(-4163)  # xlPasteValues

$ws.Range("C3").Value = "Mohaupt"

# Row 5: statement start balance line
$ws.Range("D5").Value = "KONTOSTAND AM 10.09.2023"

# Row 6: transaction 1
$ws.Range("B6").Value = "12.09."
$ws.Range("C6").Value = "13.09."
$ws.Range("D6").Value = "ABSCHLAG STROM Stadtwerke Rosenheim 24942817"
$ws.Range("E6").Value = "86,48-"

# Row 7: transaction 2
$ws.Range("B7").Value = "13.09."
$ws.Range("C7").Value = "14.09."
$ws.Range("D7").Value = "BEITRAG Allianz SE K-25485881"
$ws.Range("E7").Value = "56,76-"

# Row 8: transaction 3
$ws.Range("B8").Value = "16.09."
$ws.Range("C8").Value = "17.09."
$ws.Range("D8").Value = "KARTENZAHLUNG ARAL TANKSTELLE"
$ws.Range("E8").Value = "87,63-"

# Row 12: closing balance line
$ws.Range("D12").Value = "KONTOSTAND AM 18.09.2023"
$ws.Range("E12").Value = "230,87-"

# Row 13: next statement date
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 25.09.2023"
